$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Case 07 (row 7): the locator text was a placeholder ("true"); correct it to the real label
$ws.Range("E7").Value = "Cancelar solicitud"

# Case 08 (row 8, CP007_cancelar_solicitud): fill in the missing target-name + action cells
$ws.Range("D8").Value = "juan martin isola"
$ws.Range("E8").Value = "Agregar"
$ws.Range("E8").HorizontalAlignment = -4131
$ws.Range("E8").VerticalAlignment = -4108

# Move the view/selection down to the newly-completed row
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("E8").Select()
